# Big stimulus update:
#  - rename "face" image category to "book" (face//face_XX.jpg -> book//book_XX.jpg)
#    appears scattered across the promptFile/correctFile/dist_01File/dist_02File
#    columns (A-D)
#  - expand the abbreviated correct_ans codes (y/b/r) to full words
#    (left/center/right) in column L ("correct_ans")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Map of abbreviated answer codes -> full words (column L: "correct_ans")
$ansMap = @{ "y" = "left"; "b" = "center"; "r" = "right" }

# Stimulus image columns that can hold "face//face_NN.jpg" style paths
$imgCols = @(1, 2, 3, 4)   # A, B, C, D
$ansCol = 12                # L

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $imgCols) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            [string]$val = $v
            if ($val -like "*face//face_*") {
                $cell.Value = ($val -replace "face", "book")
            }
        }
    }

    $ansCell = $ws.Cells.Item($r, $ansCol)
    $av = $ansCell.Value2
    if ($av -ne $null) {
        [string]$aval = $av
        if ($ansMap.ContainsKey($aval)) {
            $ansCell.Value = $ansMap[$aval]
        }
    }
}
